# 1. add dropped things disappear action 2. add fly animation speed transform
# Adds a new log row (row 52) to Sheet1: date / task done / risk / duration.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A52: date text "2012.7.22" -------------------------------------------
# Typed directly, Excel/the engine auto-parses "2012.7.22" as a date serial.
# Route it through a text formula + paste-as-values so it lands as a plain
# shared string (matching how the rest of the sheet stores these "dates").
$ws.Range("A52").Formula = "=""2012.7.22"""
$ws.Range("A52").Copy() | Out-Null
$ws.Range("A52").PasteSpecial(-4163) | Out-Null

# --- B52: task description (wrapped, two lines) ----------------------------
$ws.Range("B52").Value = "加入根据飞行速度调节翅膀频率。加入落地" + [char]10 + "物品消失action"
$ws.Range("B52").WrapText = $true

# --- C52: risk note ----------------------------------------------------------
$ws.Range("C52").Value = "水晶球削球有bug。未找出"

# --- D52: work duration (hours) --------------------------------------------
$ws.Range("D52").Value = 2

# Row 52 mirrors the height of the other two-line wrapped rows.
$ws.Rows.Item(52).RowHeight = 27

# Update the active selection the way the workbook shows it after the edit.
$ws.Range("D56").Select() | Out-Null
